$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47/48: Maker and ARBITRUM swap positions (coin name + link)
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"

# Row 2
$ws.Range("D2").Value = "37.340.85"
$ws.Range("E2").Value = "  +2.36%  "

# Row 3
$ws.Range("D3").Value = "2.032.80"
$ws.Range("E3").Value = "  +4.15%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").Value = "'247.47"
$ws.Range("E5").Value = "  +1.84%  "

# Row 6
$ws.Range("D6").Value = "'0.625"
$ws.Range("E6").Value = "  -0.65%  "

# Row 7
$ws.Range("D7").Value = "'59.78"
$ws.Range("E7").Value = "  -0.10%  "

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").Value = "'0.394"
$ws.Range("E9").Value = "  +4.93%  "

# Row 10
$ws.Range("D10").Value = "'0.0808"
$ws.Range("E10").Value = "  +3.44%  "

# Row 11
$ws.Range("D11").Value = "'0.104"
$ws.Range("E11").Value = "  +0.88%  "

# Row 12
$ws.Range("D12").Value = "'15.22"
$ws.Range("E12").Value = "  +8.54%  "

# Row 13
$ws.Range("D13").Value = "2.334.79"
$ws.Range("E13").Value = "  +4.16%  "

# Row 14
$ws.Range("D14").Value = "'0.852"
$ws.Range("E14").Value = "  +1.34%  "

# Row 15
$ws.Range("D15").Value = "'21.98"
$ws.Range("E15").Value = "  +2.40%  "

# Row 16
$ws.Range("E16").Value = "  +4.18%  "

# Row 17
$ws.Range("D17").Value = "2.030.80"
$ws.Range("E17").Value = "  +3.42%  "

# Row 18
$ws.Range("D18").Value = "37.312.03"
$ws.Range("E18").Value = "  +2.56%  "

# Row 19
$ws.Range("D19").Value = "'70.44"
$ws.Range("E19").Value = "  +2.08%  "

# Row 20
$ws.Range("E20").Value = "  +1.37%  "

# Row 21
$ws.Range("E21").Value = "  +4.14%  "

# Row 22
$ws.Range("D22").Value = "'230.65"
$ws.Range("E22").Value = "  +0.75%  "

# Row 23
$ws.Range("E23").Value = "  -0.04%  "

# Row 24
$ws.Range("E24").Value = "  +5.53%  "

# Row 25
$ws.Range("D25").Value = "'2.37"
$ws.Range("E25").Value = "  +1.17%  "

# Row 26
$ws.Range("E26").Value = "  +3.66%  "

# Row 27
$ws.Range("D27").Value = "'164.28"
$ws.Range("E27").Value = "  +1.84%  "

# Row 28
$ws.Range("E28").Value = "  -4.55%  "

# Row 29
$ws.Range("D29").Value = "'19.94"
$ws.Range("E29").Value = "  +3.91%  "

# Row 30
$ws.Range("E30").Value = "  +5.15%  "

# Row 31
$ws.Range("E31").Value = "  +1.23%  "

# Row 32
$ws.Range("D32").Value = "'0.0679"
$ws.Range("E32").Value = "  +11.55%  "

# Row 33
$ws.Range("E33").Value = "  +0.63%  "

# Row 34
$ws.Range("E34").Value = "  +13.98%  "

# Row 35
$ws.Range("D35").Value = "'4.50"
$ws.Range("E35").Value = "  +1.56%  "

# Row 36
$ws.Range("D36").Value = "'3.61"
$ws.Range("E36").Value = "  +7.07%  "

# Row 37
$ws.Range("E37").Value = "  +0.08%  "

# Row 38
$ws.Range("E38").Value = "  +1.44%  "

# Row 39
$ws.Range("D39").Value = "'5.48"
$ws.Range("E39").Value = "  +1.65%  "

# Row 40
$ws.Range("D40").Value = "'2.98"
$ws.Range("E40").Value = "  +2.45%  "

# Row 41
$ws.Range("D41").Value = "'0.0975"
$ws.Range("E41").Value = "  +1.34%  "

# Row 42
$ws.Range("D42").Value = "'0.0218"
$ws.Range("E42").Value = "  +4.32%  "

# Row 43
$ws.Range("E43").Value = "  +2.48%  "

# Row 44
$ws.Range("D44").Value = "'16.71"
$ws.Range("E44").Value = "  +6.08%  "

# Row 45
$ws.Range("D45").Value = "'91.76"
$ws.Range("E45").Value = "  +4.16%  "

# Row 46
$ws.Range("D46").Value = "'7.56"
$ws.Range("E46").Value = "  +5.52%  "

# Row 47
$ws.Range("D47").Value = "'1.06"
$ws.Range("E47").Value = "  +3.87%  "

# Row 48
$ws.Range("D48").Value = "1.380.05"
$ws.Range("E48").Value = "  +1.71%  "

# Row 49
$ws.Range("E49").Value = "  +16.56%  "

# Row 50
$ws.Range("E50").Value = "  +2.45%  "

# Row 51
$ws.Range("D51").Value = "'46.20"
$ws.Range("E51").Value = "  +0.72%  "

